$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format date-text columns so values are kept as literal text
$ws.Range("Y2,AA2,Y4,AA4,Y5,AA5,Y7,AA7").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,1).Value = 81449270
$ws.Cells.Item(2,2).Value = 90082
$ws.Cells.Item(2,5).Value = 757
$ws.Cells.Item(2,7).Value = 'Hapalopilus aurantiacus'
$ws.Cells.Item(2,8).Value = '(Rostk.) Bondartsev & Singer'
$ws.Cells.Item(2,10).Value = ''
$ws.Cells.Item(2,11).Value = ''
$ws.Cells.Item(2,14).Value = ''
$ws.Cells.Item(2,17).Value = 581880.1778726063
$ws.Cells.Item(2,18).Value = 6559077.377833399
$ws.Cells.Item(2,25).Value = '2018-10-05'
$ws.Cells.Item(2,27).Value = '2018-10-05'
$ws.Cells.Item(2,29).Value = 'tall'
$ws.Cells.Item(2,32).Value = ''
$ws.Cells.Item(2,6).ClearContents()

# Row 3
$ws.Cells.Item(3,1).Value = 81447416
$ws.Cells.Item(3,2).Value = 89776
$ws.Cells.Item(3,5).Value = 6040162
$ws.Cells.Item(3,7).Value = 'Leptoporus erubescens'
$ws.Cells.Item(3,8).Value = '(Fr.) Bourdot & Galzin'
$ws.Cells.Item(3,17).Value = 581880.1778726063
$ws.Cells.Item(3,18).Value = 6559077.377833399
$ws.Cells.Item(3,6).ClearContents()

# Row 4
$ws.Cells.Item(4,1).Value = 81447718
$ws.Cells.Item(4,2).Value = 89392
$ws.Cells.Item(4,4).Value = 'NT'
$ws.Cells.Item(4,5).Value = 1202
$ws.Cells.Item(4,6).Value = 'Ullticka'
$ws.Cells.Item(4,7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(4,8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(4,17).Value = 581883.9816994623
$ws.Cells.Item(4,18).Value = 6558852.283487952
$ws.Cells.Item(4,23).Value = 'Sköldinge'
$ws.Cells.Item(4,25).Value = '2018-10-04'
$ws.Cells.Item(4,27).Value = '2018-10-04'

# Row 5
$ws.Cells.Item(5,1).Value = 81447417
$ws.Cells.Item(5,2).Value = 78098
$ws.Cells.Item(5,5).Value = 6453
$ws.Cells.Item(5,6).Value = 'Vedskivlav'
$ws.Cells.Item(5,7).Value = 'Hertelidea botryosa'
$ws.Cells.Item(5,8).Value = '(Fr.) Printzen & Kantvilas'
$ws.Cells.Item(5,17).Value = 581908.5393505754
$ws.Cells.Item(5,18).Value = 6558976.941240109
$ws.Cells.Item(5,25).Value = '2018-10-05'
$ws.Cells.Item(5,27).Value = '2018-10-05'

# Row 6
$ws.Cells.Item(6,1).Value = 81447420
$ws.Cells.Item(6,2).Value = 89794
$ws.Cells.Item(6,4).Value = 'LC'
$ws.Cells.Item(6,5).Value = 5321
$ws.Cells.Item(6,6).Value = 'Barkticka'
$ws.Cells.Item(6,7).Value = 'Rigidoporus corticola'
$ws.Cells.Item(6,8).Value = '(Fr.) Pouzar'
$ws.Cells.Item(6,17).Value = 582184.5575369275
$ws.Cells.Item(6,18).Value = 6558949.024255753
$ws.Cells.Item(6,23).Value = 'Floda'

# Row 7
$ws.Cells.Item(7,1).Value = 81447719
$ws.Cells.Item(7,2).Value = 90676
$ws.Cells.Item(7,5).Value = 5966
$ws.Cells.Item(7,6).Value = 'Motaggsvamp'
$ws.Cells.Item(7,7).Value = 'Sarcodon squamosus'
$ws.Cells.Item(7,8).Value = '(Schaeff.) Quél.'
$ws.Cells.Item(7,17).Value = 581862.0392222989
$ws.Cells.Item(7,18).Value = 6558845.144295719
$ws.Cells.Item(7,25).Value = '2018-10-04'
$ws.Cells.Item(7,27).Value = '2018-10-04'
$ws.Cells.Item(7,10).ClearContents()
$ws.Cells.Item(7,11).ClearContents()
$ws.Cells.Item(7,14).ClearContents()
$ws.Cells.Item(7,29).ClearContents()
$ws.Cells.Item(7,32).ClearContents()

# Row 8
$ws.Cells.Item(8,1).Value = 81447414
$ws.Cells.Item(8,2).Value = 90676
$ws.Cells.Item(8,5).Value = 5966
$ws.Cells.Item(8,6).Value = 'Motaggsvamp'
$ws.Cells.Item(8,7).Value = 'Sarcodon squamosus'
$ws.Cells.Item(8,8).Value = '(Schaeff.) Quél.'
$ws.Cells.Item(8,17).Value = 581965.5323911189
$ws.Cells.Item(8,18).Value = 6558883.272502526
